# Generate Report for Handback
#
# The handback for both locales (zh-cn, de-de) came back in sync with
# en-US, so refresh the status text, the "Latest Handback DateTime"
# timestamps, and clear the now-stale "handback file is not the latest"
# warning that used to live in the "Error Detail" column. The columns
# whose displayed text got longer/shorter are then re-fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Status text used to read "Ready for handoff" everywhere it's shown: the
# Overview rollup columns for each locale, and the Status column on each
# locale's own detail sheet. Now that the handback is done, update it
# everywhere.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Latest Handback DateTime stamps move forward to this handback run.
$wsZhCn.Range("K2").Value = "2016-08-25 20:48:59"
$wsDeDe.Range("K2").Value = "2016-08-25 20:49:13"

# The handback is now current, so the "handback file is not the latest"
# warning no longer applies -- clear the Error Detail cells.
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Re-fit the columns whose displayed text changed length: the Status
# column got a longer message, and the Error Detail column lost its long
# warning text (so it shrinks back toward its header width).
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsZhCn.Columns.Item(3).AutoFit()
$wsZhCn.Columns.Item(16).AutoFit()
$wsDeDe.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(16).AutoFit()

# Nudge the fitted widths to the precise values the report expects.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1
$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(16).ColumnWidth = 12.85
$wsDeDe.Columns.Item(16).ColumnWidth = 12.85
